$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns right by one.
$ws.Columns("A").Insert()

# Populate the new "Match ID" column.
$ws.Range("A1").Value2 = "Match ID"
$ws.Range("A4:A15").Value2 = 15

# Apply bold font (no border) to the header + visible data rows of the new column,
# matching the new cell style used elsewhere in the workbook.
$ws.Range("A1:A14").Font.Bold = $true

# Update the selection to match the authored state.
$ws.Range("A1:A14").Select()
